$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update GEISA's saldo (row 2) from 43169.6 to 50000
$ws.Range("C2").Value = 50000

# Delete BERTILLA's row (row 3) entirely
$ws.Rows(3).Delete()

# After deleting row 3, the former rows 8,9,10 (FABIANA, BLUEMETRIX, DAIANNE)
# shift up to rows 7,8,9. Delete from the bottom up to keep indices stable.
$ws.Rows(9).Delete()
$ws.Rows(8).Delete()
$ws.Rows(7).Delete()
